$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 89.916664
$ws.Range("I12").Value = 89
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 89
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = 81

$ws.Range("H97").Value = 2710.25
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2710.25
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 8130.75
$ws.Range("N97").Value = -9122.75

$ws.Range("H98").Value = 747.25806
$ws.Range("I98").Value = 623.62067
$ws.Range("J98").Value = 2540
$ws.Range("K98").Value = 623.62067
$ws.Range("L98").Value = 2540
$ws.Range("M98").Value = 874.37933
$ws.Range("N98").Value = -5536

$ws.Range("H114").Value = 79000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 79000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 79000
$ws.Range("N114").Value = -87678

$ws.Range("H122").Value = 747.25806
$ws.Range("I122").Value = 623.62067
$ws.Range("J122").Value = 2540
$ws.Range("K122").Value = 1870.86201
$ws.Range("L122").Value = 7620
$ws.Range("M122").Value = 579.1379899999999
$ws.Range("N122").Value = -12520

$ws.Range("H134").Value = 49771.184
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 49771.184
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 49771.184
$ws.Range("N134").Value = -59911.184

$ws.Range("H135").Value = 2086.5667
$ws.Range("I135").Value = 1855.7037
$ws.Range("J135").Value = 4164.3335
$ws.Range("K135").Value = 16701.3333
$ws.Range("L135").Value = 37479.0015
$ws.Range("M135").Value = -14166.3333

$ws.Range("H137").Value = 1534.8857
$ws.Range("I137").Value = 1234
$ws.Range("J137").Value = 2287.1
$ws.Range("K137").Value = 3702
$ws.Range("L137").Value = 6861.299999999999
$ws.Range("M137").Value = -1152
$ws.Range("N137").Value = -11961.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1743.6207
$ws.Range("I45").Value = 1402.1
$ws.Range("J45").Value = 2502.5557
$ws.Range("K45").Value = 1402.1
$ws.Range("L45").Value = 2502.5557
$ws.Range("M45").Value = -1025.1
$ws.Range("N45").Value = -3256.5557

$ws.Range("H61").Value = 6573.7144
$ws.Range("I61").Value = 6924.5
$ws.Range("J61").Value = 4469
$ws.Range("K61").Value = 6924.5
$ws.Range("L61").Value = 4469
$ws.Range("M61").Value = -6712.5
$ws.Range("N61").Value = -4893

$ws.Range("H74").Value = 2112.4119
$ws.Range("I74").Value = 2361.85
$ws.Range("J74").Value = 1756.0714
$ws.Range("K74").Value = 2361.85
$ws.Range("L74").Value = 1756.0714
$ws.Range("M74").Value = -1487.85

$ws.Range("H77").Value = 2112.4119
$ws.Range("I77").Value = 2361.85
$ws.Range("J77").Value = 1756.0714
$ws.Range("K77").Value = 11809.25
$ws.Range("L77").Value = 8780.357
$ws.Range("M77").Value = -7441.25

$ws.Range("H97").Value = 7935
$ws.Range("I97").Value = 5622
$ws.Range("J97").Value = 19500
$ws.Range("K97").Value = 5622
$ws.Range("L97").Value = 19500
$ws.Range("M97").Value = -5126

$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H131").Value = 63000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 63000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 63000
$ws.Range("N131").Value = -73080

$ws.Range("H132").Value = 5308.162
$ws.Range("I132").Value = 4312.222
$ws.Range("J132").Value = 7997.2
$ws.Range("K132").Value = 12936.666
$ws.Range("L132").Value = 23991.6
$ws.Range("M132").Value = -10406.666

$ws.Range("H136").Value = 6573.7144
$ws.Range("I136").Value = 6924.5
$ws.Range("J136").Value = 4469
$ws.Range("K136").Value = 20773.5
$ws.Range("L136").Value = 13407
$ws.Range("M136").Value = -18223.5
$ws.Range("N136").Value = -18507

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2680.842
$ws.Range("I20").Value = 2285.7778
$ws.Range("J20").Value = 3650.5454
$ws.Range("K20").Value = 2285.7778
$ws.Range("L20").Value = 3650.5454
$ws.Range("M20").Value = -2038.7778

$ws.Range("H86").Value = 1702158.2
$ws.Range("I86").Value = 2126822.8
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 2126822.8
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -2125699.8

$ws.Range("H89").Value = 1702158.2
$ws.Range("I89").Value = 2126822.8
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 10634114
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -10628498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3195.9092
$ws.Range("I58").Value = 3897
$ws.Range("J58").Value = 1326.3334
$ws.Range("K58").Value = 3897
$ws.Range("L58").Value = 1326.3334
$ws.Range("M58").Value = -3694

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H86").Value = 4449.6924
$ws.Range("I86").Value = 3699
$ws.Range("J86").Value = 5650.8
$ws.Range("K86").Value = 3699
$ws.Range("L86").Value = 5650.8
$ws.Range("M86").Value = -2576

$ws.Range("H89").Value = 4449.6924
$ws.Range("I89").Value = 3699
$ws.Range("J89").Value = 5650.8
$ws.Range("K89").Value = 18495
$ws.Range("L89").Value = 28254
$ws.Range("M89").Value = -12879

$ws.Range("H132").Value = 1733.2858
$ws.Range("I132").Value = 1493.4445
$ws.Range("J132").Value = 3172.3333
$ws.Range("K132").Value = 4480.333500000001
$ws.Range("L132").Value = 9516.999899999999
$ws.Range("M132").Value = -1950.333500000001

$ws.Range("H136").Value = 3195.9092
$ws.Range("I136").Value = 3897
$ws.Range("J136").Value = 1326.3334
$ws.Range("K136").Value = 11691
$ws.Range("L136").Value = 3979.0002
$ws.Range("M136").Value = -9141

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 450838.22
$ws.Range("I5").Value = 38657.617
$ws.Range("J5").Value = 3336102.2
$ws.Range("K5").Value = 115972.851
$ws.Range("L5").Value = 10008306.6
$ws.Range("M5").Value = -115860.851

$ws.Range("H23").Value = 2561.6
$ws.Range("I23").Value = 2074.5
$ws.Range("J23").Value = 2683.375
$ws.Range("K23").Value = 6223.5
$ws.Range("L23").Value = 8050.125
$ws.Range("M23").Value = -5988.5
$ws.Range("N23").Value = -8520.125

$ws.Range("H135").Value = 450838.22
$ws.Range("I135").Value = 38657.617
$ws.Range("J135").Value = 3336102.2
$ws.Range("K135").Value = 347918.553
$ws.Range("L135").Value = 30024919.8
$ws.Range("M135").Value = -345383.553

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 243.5
$ws.Range("I17").Value = 243.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 243.5
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -75.5
$ws.Range("N17").ClearContents()

$ws.Range("H25").Value = 3150
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 3150
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 3150
$ws.Range("N25").Value = -4208

$ws.Range("H40").Value = 38999.5
$ws.Range("I40").Value = 10000
$ws.Range("J40").Value = 67999
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 67999
$ws.Range("M40").Value = -9849
$ws.Range("N40").Value = -68301

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H47").Value = 18374.143
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 18374.143
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 18374.143
$ws.Range("N47").Value = -19510.143

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 3035.5356
$ws.Range("I122").Value = 2104.5
$ws.Range("J122").Value = 6449.3335
$ws.Range("K122").Value = 6313.5
$ws.Range("L122").Value = 19348.0005
$ws.Range("M122").Value = -3863.5

$ws.Range("H126").Value = 4030.1875
$ws.Range("I126").Value = 3997.4
$ws.Range("J126").Value = 4045.0908
$ws.Range("K126").Value = 11992.2
$ws.Range("L126").Value = 12135.2724
$ws.Range("M126").Value = -9522.200000000001
$ws.Range("N126").Value = -17075.2724

$ws.Range("H132").Value = 40178.55
$ws.Range("I132").Value = 7201.4287
$ws.Range("J132").Value = 126743.5
$ws.Range("K132").Value = 21604.2861
$ws.Range("L132").Value = 380230.5
$ws.Range("M132").Value = -19074.2861
$ws.Range("N132").Value = -385290.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 360.35715
$ws.Range("I16").Value = 330.41666
$ws.Range("J16").Value = 540
$ws.Range("K16").Value = 330.41666
$ws.Range("L16").Value = 540
$ws.Range("M16").Value = -160.41666

$ws.Range("H55").Value = 998.1579
$ws.Range("I55").Value = 110.85714
$ws.Range("J55").Value = 1515.75
$ws.Range("K55").Value = 110.85714
$ws.Range("L55").Value = 1515.75
$ws.Range("M55").Value = 62.14286
$ws.Range("N55").Value = -1861.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1325.48
$ws.Range("I107").Value = 1491.7368
$ws.Range("J107").Value = 799
$ws.Range("K107").Value = 4475.2104
$ws.Range("L107").Value = 2397
$ws.Range("M107").Value = -2555.2104
$ws.Range("N107").Value = -6237

$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H126").Value = 1183.4286
$ws.Range("I126").Value = 1188.5834
$ws.Range("J126").Value = 1152.5
$ws.Range("K126").Value = 3565.7502
$ws.Range("L126").Value = 3457.5
$ws.Range("M126").Value = -1095.7502

$ws.Range("H136").Value = 54939.848
$ws.Range("I136").Value = 13993.75
$ws.Range("J136").Value = 202345.8
$ws.Range("K136").Value = 41981.25
$ws.Range("L136").Value = 607037.3999999999
$ws.Range("M136").Value = -39431.25
$ws.Range("N136").Value = -612137.3999999999

$ws.Range("H141").Value = 58518.332
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 58518.332
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 58518.332
$ws.Range("N141").Value = -68878.33199999999
